# EPBDS-8844 Customizing output of a SpreadsheetResult.
# Rename the step/values/formula labels that used a leading "*" marker
# (e.g. "*Step1") so that the asterisk is now a trailing marker instead
# (e.g. "Step1*").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C12").Value = "Values1*"
$ws.Range("B13").Value = "Step2*"
$ws.Range("B6").Value = "Step1*"
$ws.Range("C34").Value = "Values2*"
$ws.Range("D34").Value = "Formula2*"
$ws.Range("B36").Value = "Step4*"
$ws.Range("B35").Value = "Step3*"
$ws.Range("B27").Value = "Values1*"
$ws.Range("C26").Value = "Step2*"
$ws.Range("B20").Value = "Step1*"
